# Updates "合肥-漫展信息.xlsx" gh-pages data:
#  - bump "want to go" counts (column F) on sheets "展览" and "全部类型"
#  - append a new event row (合肥·银魂主题派对only2.0) to both of those sheets

$wb = $excel.ActiveWorkbook

function Update-EventSheet {
    param($ws, $lastRow)

    # --- refresh column F (想去人数) counters for existing rows ---
    $updates = @{
        2  = 1038
        3  = 315
        4  = 1428
        5  = 8600
        6  = 74
        7  = 482
        8  = 629
        9  = 265
        10 = 149
        11 = 3462
        12 = 47
        13 = 349
        15 = 1031
        17 = 1100
        18 = 304
        19 = 178
        20 = 2165
    }
    foreach ($r in $updates.Keys) {
        $ws.Cells.Item($r, 6).Value = $updates[$r]
    }

    # --- append the new row ---
    $newRow = $lastRow + 1

    # carry over column A's bold/border/center style used for the index cells
    $ws.Range("A$lastRow").Copy() | Out-Null
    $ws.Range("A$newRow").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($newRow, 1).Value = ($newRow - 1)

    # column B holds a literal "yyyy-mm-dd" text label, not a real date --
    # force text number formatting for the assignment then drop back to the
    # sheet's normal (unstyled) look, same as the other data cells.
    $ws.Range("B$newRow").NumberFormat = "@"
    $ws.Range("B$newRow").Value = "2024-08-17"
    $ws.Range("B$newRow").Style = "Normal"

    $ws.Cells.Item($newRow, 3).Value = "合肥·银魂主题派对only2.0"
    $ws.Cells.Item($newRow, 4).Value = "长江东路1137号圣大国际商贸中心2-301室 梦田音乐LiveHouse(合肥店)"
    $ws.Cells.Item($newRow, 5).Value = "2024.08.17 13:00-08.17 18:00"
    $ws.Cells.Item($newRow, 6).Value = 0
    $ws.Cells.Item($newRow, 7).Value = 88
    $ws.Cells.Item($newRow, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87173"
    $ws.Cells.Item($newRow, 9).Value = "//i2.hdslb.com/bfs/openplatform/202406/aSc8SoTl1718078234193.png"
}

# Sheet "展览" (1st tab): rows 2-20 existing, new row 21
$wsExhibit = $wb.Worksheets.Item(1)
Update-EventSheet $wsExhibit 20

# Sheet "全部类型" (4th tab): rows 2-20 existing, new row 22 (row 21 already
# holds the unrelated "演出" entry that is untouched by this update)
$wsAll = $wb.Worksheets.Item(4)
Update-EventSheet $wsAll 21
